$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look like plain numbers to Excel; force
# them to Text format first so COM stores them as strings (matching the
# inline-string cell type in the workbook) instead of auto-converting to numbers.
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D13", "D14", "D15", "D18", "D19", "D22", "D23", "D25", "D26", "D27", "D28", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.489.02"
$ws.Range("E2").Value = "  -2.67%  "

$ws.Range("D3").Value = "1.745.28"
$ws.Range("E3").Value = "  -3.35%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "321.70"
$ws.Range("E5").Value = "  -4.16%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "0.4206"
$ws.Range("E7").Value = "  -8.75%  "

$ws.Range("E8").Value = "  -4.10%  "

$ws.Range("D9").Value = "45.47"
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("D10").Value = "0.07404"
$ws.Range("E10").Value = "  -3.16%  "

$ws.Range("E11").Value = "  -3.54%  "

$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "21.39"
$ws.Range("E13").Value = "  -4.40%  "

$ws.Range("D14").Value = "6.098"
$ws.Range("E14").Value = "  -3.83%  "

$ws.Range("D15").Value = "7.175"
$ws.Range("E15").Value = "  -3.99%  "

$ws.Range("D16").Value = "1.744.67"
$ws.Range("E16").Value = "  -3.42%  "

$ws.Range("E17").Value = "  -3.29%  "

$ws.Range("D18").Value = "87.77"
$ws.Range("E18").Value = "  +7.17%  "

$ws.Range("D19").Value = "0.06104"
$ws.Range("E19").Value = "  -9.16%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("E21").Value = "  -3.54%  "

$ws.Range("D22").Value = "6.075"
$ws.Range("E22").Value = "  -5.32%  "

$ws.Range("D23").Value = "0.5267"
$ws.Range("E23").Value = "  -5.35%  "

$ws.Range("D24").Value = "27.512.39"
$ws.Range("E24").Value = "  -2.57%  "

$ws.Range("D25").Value = "11.44"
$ws.Range("E25").Value = "  -3.52%  "

$ws.Range("D26").Value = "2.334"
$ws.Range("E26").Value = "  -3.09%  "

$ws.Range("D27").Value = "20.31"
$ws.Range("E27").Value = "  -2.40%  "

$ws.Range("D28").Value = "152.42"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").Value = "1.941.90"
$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("D31").Value = "125.58"
$ws.Range("E31").Value = "  -5.83%  "

$ws.Range("D32").Value = "1.200"
$ws.Range("E32").Value = "  -4.37%  "

$ws.Range("D33").Value = "5.653"
$ws.Range("E33").Value = "  -3.62%  "

$ws.Range("D34").Value = "0.09115"

$ws.Range("D35").Value = "3.632"
$ws.Range("E35").Value = "  -9.82%  "

$ws.Range("D36").Value = "12.57"
$ws.Range("E36").Value = "  +3.72%  "

$ws.Range("D37").Value = "0.02290"
$ws.Range("E37").Value = "  -2.64%  "

$ws.Range("D38").Value = "0.2130"
$ws.Range("E38").Value = "  -4.19%  "

$ws.Range("D39").Value = "5.064"
$ws.Range("E39").Value = "  -3.59%  "

$ws.Range("D40").Value = "0.06042"
$ws.Range("E40").Value = "  -5.05%  "

$ws.Range("D41").Value = "0.6358"
$ws.Range("E41").Value = "  -4.21%  "

$ws.Range("E42").Value = "  -3.92%  "

$ws.Range("D43").Value = "1.439"
$ws.Range("E43").Value = "  -3.99%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").Value = "7.879"
$ws.Range("E45").Value = "  -4.59%  "

$ws.Range("D46").Value = "13.61"
$ws.Range("E46").Value = "  -5.52%  "

$ws.Range("D47").Value = "3.710"
$ws.Range("E47").Value = "  -2.97%  "

$ws.Range("D48").Value = "0.5837"
$ws.Range("E48").Value = "  -4.57%  "

$ws.Range("D49").Value = "124.77"
$ws.Range("E49").Value = "  -3.77%  "

$ws.Range("E50").Value = "  -5.28%  "

$ws.Range("D51").Value = "0.06843"
$ws.Range("E51").Value = "  -4.49%  "

# Restore the original (default) cell formatting now that the text values
# are locked in, so no stray style indices are introduced.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
